$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "South Austin 1"
$ws.Range("B21").Value = "AUS P1"

$ws.Range("A22").Value = "South Austin 2"
$ws.Range("B22").Value = "AUS V4"

$ws.Range("A23").Value = "Walnut Creek"
$ws.Range("B23").Value = "AUS M9"

$ws.Range("F15").Select()
